$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "DZIENNICZEK WDZIĘCZNOŚCI" label used on the category/x axis is being
# shortened to "DZIENNICZEK" (so that axis labels fit better once their
# orientation is changed). Update the two cells that held the long label.
$ws.Range("D3").Value = "DZIENNICZEK"
$ws.Range("D6").Value = "DZIENNICZEK"

# Update the active selection to D6, matching the edited workbook.
$ws.Range("D6").Select()
